$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 10:35"

# --- Insert "Siria" row between "Trinidad yTobago" (row 168) and the old
#     "Aruba" row (row 169), pushing Aruba..Comoras down one row each, and
#     give Siria its refreshed case counts. Malaui (row 175) is untouched. ---

# Row 174 (was Siria) becomes "Comoras" with Comoras' old data (shifted from 173)
$ws.Range("A174").Value = "Comoras"
$ws.Range("B174").Value = 87
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 21
$ws.Range("E174").Value = 65
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 1

# Row 173 (was Comoras) becomes "Barbados" with Barbados' old data (shifted from 172)
$ws.Range("A173").Value = "Barbados"
$ws.Range("B173").Value = 92
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 70
$ws.Range("E173").Value = 15
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 7

# Row 172 (was Barbados) becomes "Monaco" with Monaco's old data (shifted from 171)
$ws.Range("A172").Value = "Monaco"
$ws.Range("B172").Value = 98
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 90
$ws.Range("E172").Value = 4
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 4

# Row 171 (was Monaco) becomes "Bahamas" with Bahamas' old data (shifted from 170)
$ws.Range("A171").Value = "Bahamas"
$ws.Range("B171").Value = 100
$ws.Range("C171").Value = 0
$ws.Range("D171").Value = 46
$ws.Range("E171").Value = 43
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 11

# Row 170 (was Bahamas) becomes "Aruba" with Aruba's old data (shifted from 169)
$ws.Range("A170").Value = "Aruba"
$ws.Range("B170").Value = 101
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 95
$ws.Range("E170").Value = 3
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 3

# Row 169 (was Aruba) becomes "Siria" with brand-new updated data
$ws.Range("A169").Value = "Siria"
$ws.Range("B169").Value = 106
$ws.Range("C169").Value = 20
$ws.Range("D169").Value = 41
$ws.Range("E169").Value = 61
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 4

# --- Updated case counts for several countries (data refresh) ---

# Row 13: India
$ws.Range("B13").Value = 139237
$ws.Range("C13").Value = 701
$ws.Range("D13").Value = 57745
$ws.Range("E13").Value = 77468

# Row 36: Polonia
$ws.Range("B36").Value = 21440
$ws.Range("C36").Value = 114
$ws.Range("E36").Value = 11168

# Row 40: Rumania
$ws.Range("E40").Value = 5483
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 1188

# Row 44: Austria
$ws.Range("B44").Value = 16539
$ws.Range("C44").Value = 36
$ws.Range("D44").Value = 15138
$ws.Range("E44").Value = 760
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 641

# Row 46: Filipinas
$ws.Range("B46").Value = 14319
$ws.Range("C46").Value = 284
$ws.Range("D46").Value = 3323
$ws.Range("E46").Value = 10123
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 873

# Row 103: Sri Lanka
$ws.Range("D103").Value = 695
$ws.Range("E103").Value = 437
